$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
# (matches the original workbook's text-typed Price column)
$textCells = @("D5", "D8", "D9", "D13", "D15", "D19", "D21", "D22", "D24", "D26", "D28", "D31", "D33", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '36.257.53'
$ws.Range("E2").Value = '  +1.83%  '

$ws.Range("D3").Value = '2.002.64'
$ws.Range("E3").Value = '  +5.56%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '243.60'
$ws.Range("E5").Value = '  -1.51%  '

$ws.Range("E6").Value = '  -4.78%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '44.15'
$ws.Range("E8").Value = '  +1.71%  '

$ws.Range("D9").Value = '61.65'
$ws.Range("E9").Value = '  +7.36%  '

$ws.Range("E10").Value = '  +1.16%  '

$ws.Range("E11").Value = '  -6.04%  '

$ws.Range("E12").Value = '  -0.68%  '

$ws.Range("D13").Value = '14.33'
$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("D14").Value = '2.294.03'
$ws.Range("E14").Value = '  +5.46%  '

$ws.Range("D15").Value = '0.799'
$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").Value = '1.999.43'
$ws.Range("E16").Value = '  +5.49%  '

$ws.Range("E17").Value = '  -3.37%  '

$ws.Range("D18").Value = '36.226.87'
$ws.Range("E18").Value = '  +1.62%  '

$ws.Range("D19").Value = '70.96'
$ws.Range("E19").Value = '  -3.83%  '

$ws.Range("D20").Value = '0.0₃0809'
$ws.Range("E20").Value = '  -2.82%  '

$ws.Range("D21").Value = '12.73'
$ws.Range("E21").Value = '  -2.17%  '

$ws.Range("D22").Value = '236.13'
$ws.Range("E22").Value = '  -3.97%  '

$ws.Range("E23").Value = '  -6.71%  '

$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("E25").Value = '  -9.63%  '

$ws.Range("D26").Value = '164.73'
$ws.Range("E26").Value = '  -1.29%  '

$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("D28").Value = '19.55'
$ws.Range("E28").Value = '  +6.42%  '

$ws.Range("E29").Value = '  -11.26%  '

$ws.Range("E30").Value = '  -6.38%  '

$ws.Range("D31").Value = '21.96'
$ws.Range("E31").Value = '  +61.74%  '

$ws.Range("E32").Value = '  -1.42%  '

$ws.Range("D33").Value = '0.0579'
$ws.Range("E33").Value = '  -3.88%  '

$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.0826'
$ws.Range("E36").Value = '  +12.33%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '3.94'
$ws.Range("E37").Value = '  -7.43%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '2.09'
$ws.Range("E38").Value = '  +6.47%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.845'
$ws.Range("E39").Value = '  -1.24%  '

$ws.Range("E40").Value = '  -11.17%  '

$ws.Range("D41").Value = '0.0213'
$ws.Range("E41").Value = '  -5.73%  '

$ws.Range("E42").Value = '  +1.88%  '

$ws.Range("D43").Value = '94.95'
$ws.Range("E43").Value = '  -4.56%  '

$ws.Range("D44").Value = '2.75'
$ws.Range("E44").Value = '  +15.39%  '

$ws.Range("D45").Value = '15.86'
$ws.Range("E45").Value = '  -7.57%  '

$ws.Range("D46").Value = '1.304.70'
$ws.Range("E46").Value = '  -1.30%  '

$ws.Range("D47").Value = '0.0813'
$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  +0.48%  '

$ws.Range("D49").Value = '2.182.30'
$ws.Range("E49").Value = '  +5.00%  '

$ws.Range("E50").Value = '  -8.86%  '

$ws.Range("D51").Value = '3.80'
$ws.Range("E51").Value = '  +13.03%  '
